# ---------------------------------------------------------------------------
# "made new quality ratings tab"
#
# 1) Two subjects that lived on the separate "Pilot 2 - Alums" sheet
#    (PD2_190805_03 and PD2_190806_04) are merged into "Pilot 2" itself,
#    in subID order, pushing the three rows that used to be
#    32/33/34 (PD2_190806_01/02/03) down and renumbering the running
#    "count" column sequentially.
# 2) The now-redundant "Pilot 2 - Alums" sheet is replaced by a new
#    "Quality Ratings" summary sheet: subID + drawing1_prompt +
#    drawing2_prompt for every row of the (now 36-row) "Pilot 2" sheet.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Pilot 2")

# ---------------------------------------------------------------------------
# Step 1: insert a fresh row at 33 (old rows 33-35 slide down to 34-36),
# and fill it with the first alum record (PD2_190805_03).
# ---------------------------------------------------------------------------
$ws2.Rows.Item(33).Insert()

$ws2.Cells.Item(33, 1).Value  = 32
$ws2.Cells.Item(33, 2).Value  = "PD2_190805_03"
$ws2.Cells.Item(33, 3).Value  = 43682
$ws2.Cells.Item(33, 4).Value  = 6.73
$ws2.Cells.Item(33, 5).Value  = "m"
$ws2.Cells.Item(33, 6).Value  = "ba693"
$ws2.Cells.Item(33, 7).Value  = "bing"
$ws2.Cells.Item(33, 8).Value  = "molly"
$ws2.Cells.Item(33, 9).Value  = "yes"
$ws2.Cells.Item(33, 10).Value = "good"
$ws2.Cells.Item(33, 11).Value = "because it looks like they spent more time"
$ws2.Cells.Item(33, 12).Value = "pass"
$ws2.Cells.Item(33, 13).Value = "linda"
$ws2.Cells.Item(33, 14).Value = "karen"
$ws2.Cells.Item(33, 15).Value = "linda"
$ws2.Cells.Item(33, 16).Value = "pass"
$ws2.Cells.Item(33, 17).Value = "no"
$ws2.Cells.Item(33, 18).Value = "NA"
$ws2.Cells.Item(33, 19).Value = "pass"
$ws2.Cells.Item(33, 20).Value = "no"
$ws2.Cells.Item(33, 21).Value = "NA"
$ws2.Cells.Item(33, 22).Value = "flower"
$ws2.Cells.Item(33, 23).Value = "flower"
$ws2.Cells.Item(33, 24).Value = "NA"
$ws2.Cells.Item(33, 25).Value = "flower"
$ws2.Cells.Item(33, 26).Value = "three flowers"
$ws2.Cells.Item(33, 27).Value = "NA"
$ws2.Cells.Item(33, 28).Value = "linda"
$ws2.Cells.Item(33, 29).Value = "because for karen she gave stars for every one but for teacher linda she only gave the real ones"
$ws2.Cells.Item(33, 30).Value = "linda"
$ws2.Cells.Item(33, 31).Value = "well because she wants to, because she actually, she does the real ones"
$ws2.Cells.Item(33, 32).Value = "perfect example if we need one "

# ---------------------------------------------------------------------------
# Step 2: renumber the "count" column for the three rows that slid down
# (old counts 32/33/34 -> new counts 33/34/35 at rows 34/35/36).
# ---------------------------------------------------------------------------
$ws2.Cells.Item(34, 1).Value = 33
$ws2.Cells.Item(35, 1).Value = 34
$ws2.Cells.Item(36, 1).Value = 35

# ---------------------------------------------------------------------------
# Step 3: append the second alum record (PD2_190806_04) as new row 37.
# ---------------------------------------------------------------------------
$ws2.Cells.Item(37, 1).Value  = 36
$ws2.Cells.Item(37, 2).Value  = "PD2_190806_04"

# Row 37 is a brand-new appended row (not produced by Rows.Insert), so it
# doesn't inherit the date number-format from the row above automatically;
# copy formats only from C34 (one of the existing date cells) first.
$ws2.Range("C34").Copy()
$ws2.Range("C37").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws2.Cells.Item(37, 3).Value  = 43683
$ws2.Cells.Item(37, 4).Value  = 6.51
$ws2.Cells.Item(37, 5).Value  = "f"
$ws2.Cells.Item(37, 6).Value  = "ba747"
$ws2.Cells.Item(37, 7).Value  = "bing"
$ws2.Cells.Item(37, 8).Value  = "molly"
$ws2.Cells.Item(37, 9).Value  = "yes"
$ws2.Cells.Item(37, 10).Value = "good"
$ws2.Cells.Item(37, 11).Value = "because a flower is not scribble scrabble"
$ws2.Cells.Item(37, 12).Value = "pass"
$ws2.Cells.Item(37, 13).Value = "karen"
$ws2.Cells.Item(37, 14).Value = "linda"
$ws2.Cells.Item(37, 15).Value = "linda"
$ws2.Cells.Item(37, 16).Value = "pass"
$ws2.Cells.Item(37, 17).Value = "no"
$ws2.Cells.Item(37, 18).Value = "NA"
$ws2.Cells.Item(37, 19).Value = "pass"
$ws2.Cells.Item(37, 20).Value = "no"
$ws2.Cells.Item(37, 21).Value = "NA"
$ws2.Cells.Item(37, 22).Value = "flower"
$ws2.Cells.Item(37, 23).Value = "flower"
$ws2.Cells.Item(37, 24).Value = "NA"
$ws2.Cells.Item(37, 25).Value = "flower"
$ws2.Cells.Item(37, 26).Value = "flower"
$ws2.Cells.Item(37, 27).Value = "NA"
$ws2.Cells.Item(37, 28).Value = "linda"
$ws2.Cells.Item(37, 29).Value = "because I haven't been practicing on drawing that kind of flower"
$ws2.Cells.Item(37, 30).Value = "karen"
$ws2.Cells.Item(37, 31).Value = "I don't know "

# ---------------------------------------------------------------------------
# Step 4: update the "Pilot 2" sheet view (scrolled/selected over to
# column Y, matching the edited workbook).
# ---------------------------------------------------------------------------
$ws2.Activate()
$excel.ActiveWindow.ScrollColumn = 22
[void]$ws2.Range("Y1:Y1048576").Select()

# ---------------------------------------------------------------------------
# Step 5: drop the old "Pilot 2 - Alums" sheet and create a new
# "Quality Ratings" sheet right after "Pilot 2", in the same tab slot.
# ---------------------------------------------------------------------------
$oldAlums = $wb.Worksheets.Item("Pilot 2 - Alums")
[void]$oldAlums.Delete()

$qr = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$qr.Name = "Quality Ratings"

$qr.Range("A1").Value = "subID"
$qr.Range("B1").Value = "drawing1_prompt"
$qr.Range("C1").Value = "drawing2_prompt"

for ($r = 2; $r -le 37; $r++) {
    $subid = $ws2.Cells.Item($r, 2).Value2
    $v     = $ws2.Cells.Item($r, 22).Value2
    $y     = $ws2.Cells.Item($r, 25).Value2
    $qr.Cells.Item($r, 1).Value = $subid
    $qr.Cells.Item($r, 2).Value = $v
    $qr.Cells.Item($r, 3).Value = $y
}

$qr.Columns.Item(1).ColumnWidth = 15.1640625
$qr.Columns.Item(2).ColumnWidth = 29
$qr.Columns.Item(3).ColumnWidth = 29.83203125

[void]$qr.Range("E9").Select()
